# Insert a new column before column A to hold the "ID" identifier values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Copy the header formatting from the (now shifted) B1 header cell onto the
# new A1 header cell so it keeps the same bold/centered/bordered style.
$ws.Range("B1").Copy($ws.Range("A1"))

# Set the new header text.
$ws.Range("A1").Value = "ID"

# Populate the new ID column with the sample identifiers for each row, in
# row order (rows 2-25), using an ordered array so write order is stable.
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
